$d = $word.ActiveDocument

# 1. The "_GoBack" bookmark always tracks the most recent edit location in
#    Word. Before this edit it sat at the very top of the document; once we
#    make our edit at the end of the document, Word relocates it there.
#    Remove it from its old position now (it is re-added at the new edit
#    location at the end of the script).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Locate the end of the document content (end of the final existing
#    paragraph, which ends in "...项目引入使用.") and place the insertion
#    point there.
$tail = $d.Content
$tail.Find.Execute("项目引入使用.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tail.Collapse(0)

# 3. Start a new paragraph for the additional note about unlisted .NET
#    versions, stripped of the numbered-list formatting it would otherwise
#    inherit from the preceding list paragraph.
$tail.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.ListFormat.RemoveNumbers()
$newPara.Style = $d.Styles.Item("Normal")

# 4. Type the new sentence into the freshly created paragraph.
$ip = $newPara.Range
$ip.Collapse(1)
$ip.InsertAfter("未列出的")
$ip.Collapse(0)
$ip.InsertAfter(".NET")
$ip.Collapse(0)
$ip.InsertAfter("版本，可根据需要自行更改类库工程版本使用")
$ip.Collapse(0)
$ip.InsertAfter(". .Net3.0 ")
$ip.Collapse(0)
$ip.InsertAfter("以上都可以使用")
$ip.Collapse(0)
$ip.InsertAfter(".NET3.5 ")
$ip.Collapse(0)
$ip.InsertAfter("项目更改类库版本，编译通过实现。")

# 5. Re-anchor "_GoBack" at the new end-of-document edit location.
$ip.Collapse(0)
$d.Bookmarks.Add("_GoBack", $ip)

Write-Output "done"
